# The dataset gained one new weekly observation. A new row of data is
# inserted at row 32 (pushing the existing rows 32-129 down to 33-130),
# and the new row is populated with its values; the previously last row
# (129) now also exists at row 130, fully intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32; this shifts rows 32..129 down to 33..130
# and Excel automatically extends the sheet dimension to A1:R130.
$ws.Rows.Item(32).EntireRow.Insert()

# The new row 32 needs the same "template" values the old row 32 had for the
# columns that did not change (those are now sitting in row 33, since it was
# a straight shift). Copy row 33 into row 32 to seed it.
$ws.Range("A33:R33").Copy()
$ws.Range("A32:R32").PasteSpecial()

# Now overwrite just the columns that hold the genuinely new data point.
$ws.Range("D32").Value = 44648
$ws.Range("K32").Value = 10000
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = 10000
$ws.Range("P32").Value = 3333
